$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

Write-Output ("before H3=" + $ws.Range("H3").Value())
Write-Output ("before G3=" + $ws.Range("G3").Value())
Write-Output ("before dimension-like used range=" + $ws.UsedRange.Address())
$ws.Range("H1:M29").EntireColumn.Delete()
Write-Output ("after H3=" + $ws.Range("H3").Value())
Write-Output ("after G3=" + $ws.Range("G3").Value())
Write-Output ("after used range=" + $ws.UsedRange.Address())
